$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: RPC filter - interface name changed from 1/1/c11/2 to 1/1/1
$g2 = @"
<get>
  <filter>
    <interfaces xmlns="http://openconfig.net/yang/interfaces">
      <interface>
        <name>1/1/1</name>
        <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
          <state>
            <port-speed></port-speed>
          </state>
        </ethernet>
      </interface>
    </interfaces>
  </filter>
</get>
"@
$ws.Range("G2").Value = $g2

# J2: rpc-reply - new message-id, interface name changed, and ethernet/state/port-speed added
$j2 = @"
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:6f8c9066-3c58-43a2-bd77-d7c99ceaed8e" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <data>
        <interfaces xmlns="http://openconfig.net/yang/interfaces">
            <interface>
                <name>1/1/1</name>
                <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
                    <state>
                        <port-speed>SPEED_100MB</port-speed>
                    </state>
                </ethernet>
            </interface>
        </interfaces>
    </data>
</rpc-reply>
"@
$ws.Range("J2").Value = $j2
